$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (239-244): date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @{ Row = 239; A = 44313; B = 1; C = 19; D = 175.0345462920313 },
    @{ Row = 240; A = 44314; B = 0; C = 19; D = 175.0345462920313 },
    @{ Row = 241; A = 44315; B = 4; C = 13; D = 119.7604790419162 },
    @{ Row = 242; A = 44316; B = 1; C = 11; D = 101.3357899585444 },
    @{ Row = 243; A = 44317; B = 2; C = 10; D = 92.12344541685859 },
    @{ Row = 244; A = 44318; B = 2; C = 12; D = 110.5481345002303 }
)

# Copy the formatting of the last existing row's date cell (A238) down so
# the new rows inherit the same style (date format + border); columns
# B-D keep the default (unstyled) formatting, matching the source rows.
$lastRow = 238
$ws.Range("A$lastRow").Copy() | Out-Null

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A${row}").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}

$excel.CutCopyMode = 0
